# Update de bases das ligas - Mexico Liga de Expansion (14-04-2024 18:28)
#
# Changes applied:
#   1. Rows 91 and 92 are two rows for the same match day that had their
#      data swapped (match id, teams, odds, etc.) - columns B..AC exchange
#      between the two rows while column A (row index) stays put.
#   2. Rows 186 and 187 likewise swap their B..AC data between each other.
#   3. Row 231 (the last data row) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($sheet, $row1, $row2, $firstCol, $lastCol) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $sheet.Cells.Item($row1, $c)
        $cell2 = $sheet.Cells.Item($row2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Columns B (2) through AC (29) hold the per-match data; column A (1) is the
# running index and is left untouched on both rows.
Swap-RowData $ws 91 92 2 29
Swap-RowData $ws 186 187 2 29

# Row 231 is dropped completely (sheet shrinks from A1:AC231 to A1:AC230).
$ws.Rows(231).Delete()
